$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Row 2
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Row 3
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-11 share the same relative formula (D-D_prev)*B/100
$ws.Range("G4:G11").FormulaR1C1 = "=(RC[-3]-R[-1]C[-3])*RC[-5]/100"

# Update selection to match diff
$ws.Range("E5").Select()
